# Generate Report for Handoff
# Adds two newly-handed-off source files (0fa92ec6-... and 65b368d1-...) as
# new rows across the Overview / zh-cn / de-de sheets, extending their
# tables from 3 rows (2 data rows) to 5 rows (4 data rows).

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5b0fbb0411f50b0cd04e6cd6191274324ccb1bd0/e2e/"

$guid1 = "0fa92ec6-66fb-4fd8-b7fa-ad65258486cf"
$guid2 = "65b368d1-2716-4587-9e49-e1b5d9d9d862"

$zhHash1 = "3eaf6dbce872370480de97e7233011e921de122b"
$zhHash2 = "94a5a18dd32492f293b4c7995123d3ddfdb5aadc"

$handoffDate = "2016-10-27 07:22:37"
$xliffDate = "2016-10-27 07:22:26"
$noHandback = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null

# Row 4
$wsOverview.Range("A4").Value = "$guid1.md"
$wsOverview.Range("B4").Value = "e2e\$guid1.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), ($repoBase + "$guid1.md"), $null, $null, "e2e\$guid1.md") | Out-Null
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("D4").Value = ""
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = $handoffDate
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 5
$wsOverview.Range("A5").Value = "$guid2.md"
$wsOverview.Range("B5").Value = "e2e\$guid2.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), ($repoBase + "$guid2.md"), $null, $null, "e2e\$guid2.md") | Out-Null
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("D5").Value = ""
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = $handoffDate
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null
$loZh.ListRows.Add() | Out-Null

# Row 4
$wsZh.Range("A4").Value = "$guid1.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), ($repoBase + "$guid1.md"), $null, $null, "$guid1.md") | Out-Null
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "'False"
$wsZh.Range("G4").Value = "$guid1.$zhHash1.zh-cn.xlf"
$wsZh.Range("H4").Value = $xliffDate
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = ""
$wsZh.Range("J4").Value = ""
$wsZh.Range("K4").Value = $noHandback
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "'True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "'False"
$wsZh.Range("P4").Value = ""

# Row 5
$wsZh.Range("A5").Value = "$guid2.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), ($repoBase + "$guid2.md"), $null, $null, "$guid2.md") | Out-Null
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "'False"
$wsZh.Range("G5").Value = "$guid2.$zhHash2.zh-cn.xlf"
$wsZh.Range("H5").Value = $xliffDate
$wsZh.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I5").Value = ""
$wsZh.Range("J5").Value = ""
$wsZh.Range("K5").Value = $noHandback
$wsZh.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("L5").Value = ""
$wsZh.Range("M5").Value = "'True"
$wsZh.Range("N5").Value = ""
$wsZh.Range("O5").Value = "'False"
$wsZh.Range("P5").Value = ""

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null
$loDe.ListRows.Add() | Out-Null

# Row 4
$wsDe.Range("A4").Value = "$guid1.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), ($repoBase + "$guid1.md"), $null, $null, "$guid1.md") | Out-Null
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "'False"
$wsDe.Range("G4").Value = "$guid1.$zhHash1.de-de.xlf"
$wsDe.Range("H4").Value = $handoffDate
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = ""
$wsDe.Range("J4").Value = ""
$wsDe.Range("K4").Value = $noHandback
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "'True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "'False"
$wsDe.Range("P4").Value = ""

# Row 5
$wsDe.Range("A5").Value = "$guid2.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), ($repoBase + "$guid2.md"), $null, $null, "$guid2.md") | Out-Null
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "'False"
$wsDe.Range("G5").Value = "$guid2.$zhHash2.de-de.xlf"
$wsDe.Range("H5").Value = $handoffDate
$wsDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I5").Value = ""
$wsDe.Range("J5").Value = ""
$wsDe.Range("K5").Value = $noHandback
$wsDe.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("L5").Value = ""
$wsDe.Range("M5").Value = "'True"
$wsDe.Range("N5").Value = ""
$wsDe.Range("O5").Value = "'False"
$wsDe.Range("P5").Value = ""

# ---------------------------------------------------------------------
# Column width adjustments (widened to fit longer datetime values)
# ---------------------------------------------------------------------
$wsOverview.Range("E1").EntireColumn.ColumnWidth = 17.2159881591797
$wsOverview.Range("F1").EntireColumn.ColumnWidth = 17.2159881591797
$wsZh.Range("C1").EntireColumn.ColumnWidth = 17.2159881591797
$wsDe.Range("C1").EntireColumn.ColumnWidth = 17.2159881591797
